# Update EUR->ARS rate: 2025-09-05T15:19:32Z
# Append a new row of data (Fecha / Hora / Cotizacion) to the quote log.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.Cells.Item(4, 1).Row + 1

# Leading apostrophe forces these to be stored as plain text (not
# auto-converted to date/number serials) while keeping the cell's
# number format at the default ("General"/style 0).
$ws.Cells.Item($newRow, 1).Value = "'2025-09-05"
$ws.Cells.Item($newRow, 1).Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = "'15:19:32"
$ws.Cells.Item($newRow, 2).Style = "Normal"

$ws.Cells.Item($newRow, 3).Value = "1.00 EUR = 1595.8598 ARS"
